$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 'b'
$ws.Range("J4").Value = 'Acknowledge (Backchannel)'
$ws.Range("I26").Value = 'sv'
$ws.Range("J26").Value = 'Statement-opinion'
$ws.Range("I48").Value = 'sv'
$ws.Range("J48").Value = 'Statement-opinion'
$ws.Range("I49").Value = 'sv'
$ws.Range("J49").Value = 'Statement-opinion'
$ws.Range("I52").Value = 'sv'
$ws.Range("J52").Value = 'Statement-opinion'
$ws.Range("I54").Value = 'aa'
$ws.Range("J54").Value = 'Agree/Accept'
$ws.Range("I66").Value = 'sd'
$ws.Range("J66").Value = 'Statement-non-opinion'
$ws.Range("I72").Value = 'sd'
$ws.Range("J72").Value = 'Statement-non-opinion'
$ws.Range("I74").Value = 'sd'
$ws.Range("J74").Value = 'Statement-non-opinion'
$ws.Range("I92").Value = 'sv'
$ws.Range("J92").Value = 'Statement-opinion'
$ws.Range("I93").Value = 'sd'
$ws.Range("J93").Value = 'Statement-non-opinion'
$ws.Range("I99").Value = 'sd'
$ws.Range("J99").Value = 'Statement-non-opinion'
$ws.Range("I103").Value = 'aa'
$ws.Range("J103").Value = 'Agree/Accept'
$ws.Range("I131").Value = 'sd'
$ws.Range("J131").Value = 'Statement-non-opinion'
$ws.Range("I144").Value = 'sv'
$ws.Range("J144").Value = 'Statement-opinion'
$ws.Range("I150").Value = 'sd'
$ws.Range("J150").Value = 'Statement-non-opinion'
$ws.Range("I152").Value = 'sd'
$ws.Range("J152").Value = 'Statement-non-opinion'
$ws.Range("I154").Value = 'sd'
$ws.Range("J154").Value = 'Statement-non-opinion'
$ws.Range("I157").Value = 'sv'
$ws.Range("J157").Value = 'Statement-opinion'
$ws.Range("I164").Value = 'sd'
$ws.Range("J164").Value = 'Statement-non-opinion'
$ws.Range("I172").Value = 'sd'
$ws.Range("J172").Value = 'Statement-non-opinion'
$ws.Range("I173").Value = 'sv'
$ws.Range("J173").Value = 'Statement-opinion'
$ws.Range("I177").Value = 'b'
$ws.Range("J177").Value = 'Acknowledge (Backchannel)'
$ws.Range("I193").Value = 'sd'
$ws.Range("J193").Value = 'Statement-non-opinion'
$ws.Range("I200").Value = 'sd'
$ws.Range("J200").Value = 'Statement-non-opinion'
$ws.Range("I202").Value = 'sv'
$ws.Range("J202").Value = 'Statement-opinion'
$ws.Range("I208").Value = 'sd'
$ws.Range("J208").Value = 'Statement-non-opinion'
$ws.Range("I209").Value = 'sd'
$ws.Range("J209").Value = 'Statement-non-opinion'
$ws.Range("I215").Value = 'sd'
$ws.Range("J215").Value = 'Statement-non-opinion'
$ws.Range("I220").Value = 'sd'
$ws.Range("J220").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'aa'
$ws.Range("J236").Value = 'Agree/Accept'
$ws.Range("I241").Value = 'sv'
$ws.Range("J241").Value = 'Statement-opinion'
$ws.Range("I242").Value = 'sd'
$ws.Range("J242").Value = 'Statement-non-opinion'
$ws.Range("I247").Value = 'aa'
$ws.Range("J247").Value = 'Agree/Accept'
$ws.Range("I255").Value = 'sd'
$ws.Range("J255").Value = 'Statement-non-opinion'
$ws.Range("I258").Value = '%'
$ws.Range("J258").Value = 'Uninterpretable'
$ws.Range("I262").Value = 'sd'
$ws.Range("J262").Value = 'Statement-non-opinion'
$ws.Range("I275").Value = 'sv'
$ws.Range("J275").Value = 'Statement-opinion'
$ws.Range("I294").Value = 'sd'
$ws.Range("J294").Value = 'Statement-non-opinion'
$ws.Range("I304").Value = 'aa'
$ws.Range("J304").Value = 'Agree/Accept'
$ws.Range("I308").Value = 'sd'
$ws.Range("J308").Value = 'Statement-non-opinion'
$ws.Range("I320").Value = 'sd'
$ws.Range("J320").Value = 'Statement-non-opinion'
$ws.Range("I324").Value = 'aa'
$ws.Range("J324").Value = 'Agree/Accept'
$ws.Range("I325").Value = 'aa'
$ws.Range("J325").Value = 'Agree/Accept'
$ws.Range("I340").Value = 'b'
$ws.Range("J340").Value = 'Acknowledge (Backchannel)'
$ws.Range("I356").Value = 'aa'
$ws.Range("J356").Value = 'Agree/Accept'
$ws.Range("I382").Value = 'sd'
$ws.Range("J382").Value = 'Statement-non-opinion'
$ws.Range("I405").Value = 'sv'
$ws.Range("J405").Value = 'Statement-opinion'
$ws.Range("I407").Value = 'sd'
$ws.Range("J407").Value = 'Statement-non-opinion'
$ws.Range("I408").Value = 'b'
$ws.Range("J408").Value = 'Acknowledge (Backchannel)'
$ws.Range("I415").Value = 'aa'
$ws.Range("J415").Value = 'Agree/Accept'
$ws.Range("I426").Value = 'aa'
$ws.Range("J426").Value = 'Agree/Accept'
$ws.Range("I452").Value = 'sd'
$ws.Range("J452").Value = 'Statement-non-opinion'
$ws.Range("I463").Value = 'aa'
$ws.Range("J463").Value = 'Agree/Accept'
$ws.Range("I480").Value = '%'
$ws.Range("J480").Value = 'Uninterpretable'
$ws.Range("I482").Value = '%'
$ws.Range("J482").Value = 'Uninterpretable'
$ws.Range("I491").Value = 'sd'
$ws.Range("J491").Value = 'Statement-non-opinion'
$ws.Range("I495").Value = '%'
$ws.Range("J495").Value = 'Uninterpretable'
$ws.Range("I501").Value = 'aa'
$ws.Range("J501").Value = 'Agree/Accept'
$ws.Range("I508").Value = 'sd'
$ws.Range("J508").Value = 'Statement-non-opinion'
$ws.Range("I537").Value = 'sv'
$ws.Range("J537").Value = 'Statement-opinion'
$ws.Range("I539").Value = 'aa'
$ws.Range("J539").Value = 'Agree/Accept'
$ws.Range("I543").Value = 'b'
$ws.Range("J543").Value = 'Acknowledge (Backchannel)'
$ws.Range("I545").Value = '%'
$ws.Range("J545").Value = 'Uninterpretable'
$ws.Range("I558").Value = 'sd'
$ws.Range("J558").Value = 'Statement-non-opinion'
$ws.Range("I580").Value = 'sd'
$ws.Range("J580").Value = 'Statement-non-opinion'
$ws.Range("I587").Value = 'sv'
$ws.Range("J587").Value = 'Statement-opinion'
$ws.Range("I588").Value = 'sd'
$ws.Range("J588").Value = 'Statement-non-opinion'
$ws.Range("I590").Value = 'sv'
$ws.Range("J590").Value = 'Statement-opinion'
$ws.Range("I598").Value = 'sd'
$ws.Range("J598").Value = 'Statement-non-opinion'
$ws.Range("I599").Value = 'ba'
$ws.Range("J599").Value = 'Appreciation'
$ws.Range("I607").Value = 'sv'
$ws.Range("J607").Value = 'Statement-opinion'
$ws.Range("I608").Value = '%'
$ws.Range("J608").Value = 'Uninterpretable'
